# Daily refresh of the "remaining days" tracker.
#
# For every data row: E (剩余/remaining) counts down by 1 per day that
# passes, tracked against F (开始时间/start date, stored as a yyyyMMdd
# integer) for a cycle length of D (总天/total days) days. When the
# countdown would reach 0 (i.e. old E == 1), the cycle restarts: E resets
# to the full D-day count and F advances by D days (a new bucket-swap
# cycle begins).
#
# Rows whose start-date stamp isn't a well-formed yyyyMMdd value are left
# completely untouched (defensive skip for bad/legacy data entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $totalDays = $dCell.Value2
    $remaining = $eCell.Value2
    $startDate = $fCell.Value2

    if ($totalDays -eq $null -or $remaining -eq $null -or $startDate -eq $null) {
        continue
    }

    # Validate the start-date stamp; skip the row entirely if it doesn't
    # parse as a clean yyyyMMdd date (defends against corrupt entries).
    try {
        $dt = [DateTime]::ParseExact("$startDate", "yyyyMMdd", $null)
    } catch {
        continue
    }

    if ($remaining -gt 1) {
        $eCell.Value = $remaining - 1
    } elseif ($remaining -eq 1) {
        $eCell.Value = $totalDays

        $dt = $dt.AddDays($totalDays)
        $newStart = ($dt.Year * 10000) + ($dt.Month * 100) + $dt.Day
        $fCell.Value = $newStart
    }
}
